# Updated cryptos list - refresh Price (column D) and Volume(1h) (column E)
# values for the rows that changed in this data pull.
#
# Note: several "Price" values are short decimal-looking strings
# (e.g. "0.998", "214.58"). Those columns are plain text in the workbook
# (the cells already hold text such as "26.834.64" which isn't a valid
# number), so a bare numeric-looking literal would otherwise be
# auto-converted to a real number by Excel and lose its original text
# representation. Prefixing with a leading apostrophe forces Excel to
# keep/store the value as literal text, matching the original formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.867.35'
$ws.Range("E2").Value = '  +2.19%  '
$ws.Range("D3").Value = '1.639.99'
$ws.Range("E3").Value = '  +2.41%  '
$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = '  -0.26%  '
$ws.Range("D5").Value = "'214.58"
$ws.Range("E5").Value = '  +0.97%  '
$ws.Range("E6").Value = '  +0.45%  '
$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = '  -0.35%  '
$ws.Range("E8").Value = '  +1.68%  '
$ws.Range("E9").Value = '  +0.95%  '
$ws.Range("D10").Value = "'19.48"
$ws.Range("E10").Value = '  +2.87%  '
$ws.Range("E11").Value = '  +0.37%  '
$ws.Range("D12").Value = '1.868.43'
$ws.Range("E12").Value = '  +2.24%  '
$ws.Range("D13").Value = '1.638.75'
$ws.Range("E13").Value = '  +2.24%  '
$ws.Range("E14").Value = '  +1.28%  '
$ws.Range("E15").Value = '  +2.00%  '
$ws.Range("D16").Value = "'64.63"
$ws.Range("E16").Value = '  +1.57%  '
$ws.Range("D17").Value = "'241.98"
$ws.Range("E17").Value = '  +6.24%  '
$ws.Range("D18").Value = '26.846.74'
$ws.Range("E18").Value = '  +2.07%  '
$ws.Range("D19").Value = "'7.88"
$ws.Range("E19").Value = '  +3.39%  '
$ws.Range("E20").Value = '  +0.74%  '
$ws.Range("E21").Value = '  -0.15%  '
$ws.Range("D22").Value = "'4.38"
$ws.Range("E22").Value = '  +1.68%  '
$ws.Range("D23").Value = "'2.27"
$ws.Range("E23").Value = '  +4.49%  '
$ws.Range("D24").Value = "'9.24"
$ws.Range("E24").Value = '  +3.16%  '
$ws.Range("D25").Value = "'145.83"
$ws.Range("E25").Value = '  +0.25%  '
$ws.Range("E26").Value = '  -0.25%  '
$ws.Range("D27").Value = "'7.08"
$ws.Range("E27").Value = '  +1.94%  '
$ws.Range("E28").Value = '  +0.75%  '
$ws.Range("E29").Value = '  +2.38%  '
$ws.Range("E30").Value = '  +0.10%  '
$ws.Range("E31").Value = '  +0.41%  '
$ws.Range("E32").Value = '  +2.05%  '
$ws.Range("D33").Value = '1.504.99'
$ws.Range("E33").Value = '  +4.26%  '
$ws.Range("D34").Value = "'3.04"
$ws.Range("E34").Value = '  +2.44%  '
$ws.Range("E35").Value = '  +6.63%  '
$ws.Range("E36").Value = '  -0.44%  '
$ws.Range("D37").Value = "'0.575"
$ws.Range("E37").Value = '  +1.74%  '
$ws.Range("E38").Value = '  +1.97%  '
$ws.Range("D39").Value = "'0.863"
$ws.Range("E39").Value = '  +4.78%  '
$ws.Range("D40").Value = "'5.97"
$ws.Range("E40").Value = '  +2.71%  '
$ws.Range("E41").Value = '  -0.28%  '
$ws.Range("E42").Value = '  +1.56%  '
$ws.Range("D43").Value = "'64.16"
$ws.Range("E43").Value = '  +5.47%  '
$ws.Range("D44").Value = '1.779.15'
$ws.Range("E44").Value = '  +2.13%  '
$ws.Range("D45").Value = "'0.767"
$ws.Range("E45").Value = '  +1.03%  '
$ws.Range("E46").Value = '  -0.65%  '
$ws.Range("D47").Value = "'90.46"
$ws.Range("E47").Value = '  +3.16%  '
$ws.Range("E48").Value = '  +3.27%  '
$ws.Range("D49").Value = "'0.0976"
$ws.Range("E49").Value = '  +2.50%  '
$ws.Range("E50").Value = '  +0.29%  '
$ws.Range("D51").Value = "'7.49"
$ws.Range("E51").Value = '  +1.26%  '
